$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# The source workbook deduplicates identical cell text into shared
# strings. Several of the values being changed here happen to be shared
# (byte-for-byte identical) between the row for
# 2d2ba3a8-48bb-4701-952a-bd99bd6a0f60.md (row 3) and the row for
# b4a59158-700a-4cf1-a37c-811a9daa2998.md (row 5) on the same sheet, and
# in one case between the Overview sheet and the de-de sheet too. To
# reproduce a pure in-place edit of those shared-string entries (rather
# than forking off new, separate strings), every cell that currently
# shares the old text must be updated to the new text together.

# "2016-09-05 00:19:02" -> "2016-09-05 00:19:58"
# shared by Overview!G3, Overview!G5, de-de!H3, de-de!H5
$wsOverview.Range("G3").Value = "2016-09-05 00:19:58"
$wsOverview.Range("G5").Value = "2016-09-05 00:19:58"
$wsDeDe.Range("H3").Value = "2016-09-05 00:19:58"
$wsDeDe.Range("H5").Value = "2016-09-05 00:19:58"

# "ht" -> "mt"
# shared by zh-cn!E3, zh-cn!E5, de-de!E3, de-de!E5
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# "2016-09-05 00:18:57" -> "2016-09-05 00:19:53"
# shared by zh-cn!H3, zh-cn!H5
$wsZhCn.Range("H3").Value = "2016-09-05 00:19:53"
$wsZhCn.Range("H5").Value = "2016-09-05 00:19:53"

# "2016-09-05 00:19:27" -> "2016-09-05 00:20:20"
# shared by zh-cn!K3, zh-cn!K5
$wsZhCn.Range("K3").Value = "2016-09-05 00:20:20"
$wsZhCn.Range("K5").Value = "2016-09-05 00:20:20"

# "2016-09-05 00:19:34" -> "2016-09-05 00:20:28"
# shared by de-de!K3, de-de!K5
$wsDeDe.Range("K3").Value = "2016-09-05 00:20:28"
$wsDeDe.Range("K5").Value = "2016-09-05 00:20:28"
